# This workbook holds a weekly feed of fruit/vegetable market prices
# (Betarraga - Vega Central Mapocho de Santiago). A new week's record was
# added at the top of the existing price history (row 427), pushing every
# following record down by one row. The last existing record therefore now
# appears twice further down (it lands on the newly created last row, 528).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row right above the current row 427. This shifts rows
# 427-527 down to 428-528, automatically carrying their values/formatting
# along with them.
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new weekly record.
$ws.Cells.Item(427, 1).Value2  = 9
$ws.Cells.Item(427, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(427, 3).Value2  = "Metropolitana"
$ws.Cells.Item(427, 4).Value2  = 44711
$ws.Cells.Item(427, 5).Value2  = 13
$ws.Cells.Item(427, 6).Value2  = 100114014
$ws.Cells.Item(427, 7).Value2  = "Betarraga"
$ws.Cells.Item(427, 8).Value2  = "Sin especificar"
$ws.Cells.Item(427, 9).Value2  = "Primera"
$ws.Cells.Item(427, 10).Value2 = 5200
$ws.Cells.Item(427, 11).Value2 = 110
$ws.Cells.Item(427, 12).Value2 = 120
$ws.Cells.Item(427, 13).Value2 = 115
$ws.Cells.Item(427, 14).Value2 = "$/unidad"
$ws.Cells.Item(427, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(427, 16).Value2 = 115
$ws.Cells.Item(427, 17).Value2 = 1
$ws.Cells.Item(427, 18).Value2 = "Hortaliza"

# Ensure the date column keeps its date/time number format after the insert.
$ws.Cells.Item(427, 4).NumberFormat = $ws.Cells.Item(428, 4).NumberFormat
